$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folder Inventory")

# Shift rows 2..64 down to rows 3..65 (process bottom-up so sources aren't
# overwritten before they're read). Row 65's prior content is discarded and
# row 1 (header) / rows 66-75 are left untouched.
for ($r = 64; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$dst").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$dst").Value2 = $ws.Range("C$r").Value2
    $ws.Range("D$dst").Value2 = $ws.Range("D$r").Value2
    $ws.Range("E$dst").Value2 = $ws.Range("E$r").Value2
}

# New entry inserted at the top of the inventory.
$ws.Range("A2").Value2 = "Create and Publish PowerBI Dashboards & Reports"
$ws.Range("B2").Value2 = "Create and Publish PowerBI Dashboards & Reports"
$ws.Range("C2").Value2 = "2025-06-12 20:05:46 +0530"
$ws.Range("D2").Value2 = 1
$ws.Range("E2").Value2 = "Root"

# Metadata sheet: refresh generation timestamp and workflow run counter.
# (Workflow Run is stored as text, so prefix with an apostrophe to keep it
# from being auto-converted to a number, matching the source sheet's type.)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value2 = "2025-06-12 15:25:18 UTC"
$meta.Range("B5").Value2 = "'17"

# Summary sheet: most recent update now reflects the new top entry.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value2 = "2025-06-12 20:05:46 +0530"
